$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-10 (A:T). The underlying natmi computation was
# redone to also include "ECs" as a sending cluster (previously only FAPs/sCs sent
# to ECs/FAPs/sCs); this expands the ligand-receptor pair table for Bmp2-Acvr2a to
# the complete 3x3 sending x target grid (ECs, FAPs, sCs) and recomputes every
# specificity/weight column accordingly.
$data = @(
    @("ECs","Bmp2","Acvr2a","ECs",2.0,0.6666666666666666,1.445484,4.336452,0.1286708197254238,0.1286708197254238,3.0,1.0,14.61878266666667,43.856348,0.2662829816142094,0.2662829816142094,21.13121644414399,190.180947997296,0.03426284952323029,0.03426284952323029),
    @("ECs","Bmp2","Acvr2a","FAPs",2.0,0.6666666666666666,1.445484,4.336452,0.1286708197254238,0.1286708197254238,3.0,1.0,27.084169,81.25250700000001,0.4933415757187404,0.4933415757187404,39.149732942796,352.347596485164,0.06347866495236257,0.06347866495236257),
    @("ECs","Bmp2","Acvr2a","sCs",2.0,0.6666666666666666,1.445484,4.336452,0.1286708197254238,0.1286708197254238,3.0,1.0,13.19647366666667,39.589421,0.2403754426670501,0.2403754426670501,19.075291541588,171.677623874292,0.03092930524983096,0.03092930524983096),
    @("FAPs","Bmp2","Acvr2a","ECs",3.0,1.0,6.292313,18.876939,0.5601148623429528,0.5601148623429528,3.0,1.0,14.61878266666667,43.856348,0.2662829816142094,0.2662829816142094,91.98595621764133,827.8736059587719,0.149149055591114,0.149149055591114),
    @("FAPs","Bmp2","Acvr2a","FAPs",3.0,1.0,6.292313,18.876939,0.5601148623429528,0.5601148623429528,3.0,1.0,27.084169,81.25250700000001,0.4933415757187404,0.4933415757187404,170.422068692897,1533.798618236073,0.2763279487717577,0.2763279487717577),
    @("FAPs","Bmp2","Acvr2a","sCs",3.0,1.0,6.292313,18.876939,0.5601148623429528,0.5601148623429528,3.0,1.0,13.19647366666667,39.589421,0.2403754426670501,0.2403754426670501,83.03634280692434,747.3270852623191,0.1346378579800811,0.1346378579800811),
    @("sCs","Bmp2","Acvr2a","ECs",3.0,1.0,3.496172,10.488516,0.3112143179316233,0.3112143179316232,3.0,1.0,14.61878266666667,43.856348,0.2662829816142094,0.2662829816142094,51.10977863328533,459.988007699568,0.08287107649986518,0.08287107649986516),
    @("sCs","Bmp2","Acvr2a","FAPs",3.0,1.0,3.496172,10.488516,0.3112143179316233,0.3112143179316232,3.0,1.0,27.084169,81.25250700000001,0.4933415757187404,0.4933415757187404,94.69091330106801,852.2182197096122,0.1535349619946201,0.1535349619946201),
    @("sCs","Bmp2","Acvr2a","sCs",3.0,1.0,3.496172,10.488516,0.3112143179316233,0.3112143179316232,3.0,1.0,13.19647366666667,39.589421,0.2403754426670501,0.2403754426670501,46.13714173213734,415.2342755892361,0.07480827943713803,0.07480827943713801)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $colNum = $j + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}

Write-Host "Wrote $($data.Count) rows"
